$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Make room: insert 4 rows before row 10 (shifts the second table from
#    rows 10-13 down to rows 14-17). We insert at 9:12, then collapse the
#    still-empty row 12 stub so only rows 9-11 remain as "to be filled".
# ---------------------------------------------------------------------------
$ws.Rows("9:12").Insert(-4121, 0)
$ws.Rows("12:12").Delete(-4121)
$ws.Rows("12:12").Insert(-4121, 0)

# ---------------------------------------------------------------------------
# 2) Row 11 (new "Y7" / Date input row) must inherit the "last row" border
#    styling that row 8 currently still has (border-left+bottom / border-
#    right+bottom), because row 8 is about to stop being the last row of
#    the little X1/X2 table. Copy formats from row 8 BEFORE changing row 8.
# ---------------------------------------------------------------------------
$ws.Range("B8:D8").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Now fix up row 8 so it becomes a normal "middle" row (copy formats
#    from row 7, which already has the correct middle-row borders).
# ---------------------------------------------------------------------------
$ws.Range("B7:D7").Copy()
$ws.Range("B8:D8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Row 9 ("Y5", Double input) - copy formats from row 7 (middle row) and
#    from the quote-prefixed "=" style used by C6 for the C column.
# ---------------------------------------------------------------------------
$ws.Range("B7:D7").Copy()
$ws.Range("B9:D9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C6").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Row 10 ("Y6", Boolean input) - same middle-row borders for B/D; C10
#    keeps the plain default/general style (no border).
# ---------------------------------------------------------------------------
$ws.Range("B7:D7").Copy()
$ws.Range("B10:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C10").ClearFormats()

# ---------------------------------------------------------------------------
# Populate the VALUES first (the leading "'=" trick used for the literal
# "=$X1.getClass().getSimpleName()" text auto-applies a quote-prefix style,
# so we must overwrite the style with a format-paste AFTER setting values).
# ---------------------------------------------------------------------------
$ws.Range("B9").Value2 = "Y5"
$ws.Range("C9").Value2 = 1.2
$ws.Range("D9").Value2 = '''=$X1.getClass().getSimpleName()'

$ws.Range("B10").Value2 = "Y6"
$ws.Range("C10").Value2 = $true
$ws.Range("D10").Value2 = '''=$X1.getClass().getSimpleName()'

$ws.Range("B11").Value2 = "Y7"
$ws.Range("C11").Value2 = 34284
$ws.Range("D11").Value2 = '''=$X1.getClass().getSimpleName()'

# Re-apply the intended (non quote-prefixed) borders to the D column cells
# now that their literal text value is set.
$ws.Range("D7").Copy()
$ws.Range("D9:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D8").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# C11 gets the old C8 (pre-fix) bottom border plus a date number format.
$ws.Range("C11").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 6) Extend the second table's header row (now row 15) with 3 more columns
#    (F,G,H) for Y5/Y6/Y7. Column H must take over the "last column" border
#    that column E currently has; capture it before changing E.
# ---------------------------------------------------------------------------
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D15").Copy()
$ws.Range("E15:G15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F15").Value2 = '_res_.$X2$Y5'
$ws.Range("G15").Value2 = '_res_.$X2$Y6'
$ws.Range("H15").Value2 = '_res_.$X2$Y7'

# ---------------------------------------------------------------------------
# 7) Extend the "Result" row (now row 16) with F16:H16.
# ---------------------------------------------------------------------------
$ws.Range("E16").Copy()
$ws.Range("F16:H16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F16").Value2 = "Result"
$ws.Range("G16").Value2 = "Result"
$ws.Range("H16").Value2 = "Result"

# ---------------------------------------------------------------------------
# 8) Extend the type row (now row 17): B17 becomes "Integer" (was "Double"),
#    the "Double" value now lives in the new F17 cell; G17/H17 are new
#    Boolean/Date columns.
# ---------------------------------------------------------------------------
$ws.Range("E17").Copy()
$ws.Range("F17:H17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B17").Value2 = "Integer"
$ws.Range("F17").Value2 = "Double"
$ws.Range("G17").Value2 = "Boolean"
$ws.Range("H17").Value2 = "Date"

# ---------------------------------------------------------------------------
# 9) Widen the merged title cell on row 14 from B:E to B:H.
# ---------------------------------------------------------------------------
$ws.Range("B14:E14").UnMerge()
$ws.Range("B14:H14").Merge()

# ---------------------------------------------------------------------------
# 10) Column widths.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 11.33203125
$ws.Columns("D").ColumnWidth = 26.1640625
$ws.Columns("E").ColumnWidth = 11.33203125
$ws.Columns("F").ColumnWidth = 11.33203125
$ws.Columns("G").ColumnWidth = 11.33203125
$ws.Columns("H").ColumnWidth = 11.33203125
$ws.Columns("K").ColumnWidth = 18.5

# ---------------------------------------------------------------------------
# 11) Selection.
# ---------------------------------------------------------------------------
$ws.Range("K8").Select()
